# Add a new metadata row (row 5) to the "Metadata Report" sheet,
# duplicating the December 2024 / Akurana entry already present in row 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata Report")

$sourceRow = 4
$targetRow = 5
$lastCol = 41   # column AO

$ws.Cells.Item($targetRow, 1).Value = 2024
$ws.Cells.Item($targetRow, 2).Value = "DEC"
$ws.Cells.Item($targetRow, 3).Value = "31/12-01/12"

for ($col = 4; $col -le $lastCol; $col++) {
    $ws.Cells.Item($targetRow, $col).Value = $ws.Cells.Item($sourceRow, $col).Value2
}
